$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AD1:AF1, copying the
# existing header style (bold font + border + center alignment) from AC1
# so the new headers match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in season record (Wins/Losses/Ties) for every data row (2-60).
for ($row = 2; $row -le 60; $row++) {
    $ws.Cells.Item($row, 30).Value = 76
    $ws.Cells.Item($row, 31).Value = 86
    $ws.Cells.Item($row, 32).Value = 0
}
